$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Simple single-value cell text replacements (row index 1-based -> new text)
function Set-CellText($rowIndex, $text) {
    $cell = $tbl.Cell($rowIndex, 1)
    $rng = $cell.Range
    # Trim trailing cell-mark/paragraph-mark characters Word appends
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "817"

Set-CellText 6 "0.00063"
Set-CellText 7 "0.00020"
Set-CellText 8 "0.00006"
Set-CellText 9 "0.00026"
Set-CellText 10 "0.00035"
Set-CellText 11 "0.00045"
Set-CellText 12 "0.15999"

# Rows 44-46 (1-based) collapse multi-run/tab content down to a single value
Set-CellText 44 "99.88"
Set-CellText 45 "0.16"
Set-CellText 46 "130"
